# Update loading_percent values for rows 2-25 (data index 0-23),
# columns B,C,D,E,G,H,I,L,N (F,J,K,M,O remain 0 and A is unchanged)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colLetters = @("B","C","D","E","G","H","I","L","N")

$data = @(
    @(19.30175218868783, 9.359008950556611, 6.000743325601204, 9.917797993583333, 43.38720887879066, 17.4197638208352, 26.10181466441462, 10.2742648893817, 18.16255530900539),
    @(18.77987472425733, 8.783565395150836, 5.885786259652837, 9.94117660039621, 42.99118484504491, 17.42482376637523, 26.13481833111841, 10.25119504932647, 18.23688922098551),
    @(18.45693006997931, 8.409137042208412, 5.816000549909685, 9.956490195497613, 42.76296120058336, 17.43236153732904, 26.16215218737007, 10.23914361927508, 18.28447527119237),
    @(18.32491096349692, 8.251260753187422, 5.787805275672643, 9.962972182862686, 42.67380653654261, 17.4365446151194, 26.17506086717872, 10.23476699346311, 18.30435784021617),
    @(18.30297060188164, 8.224726298765942, 5.783139380338438, 9.964063114635879, 42.6592374158035, 17.43730624172073, 26.17731099999035, 10.23407261356105, 18.30768902706838),
    @(18.45515100975506, 8.407029279409469, 5.815619258992571, 9.956576635163374, 42.76174313302536, 17.4324134561506, 26.16231912299762, 10.23908242712794, 18.28474142457159),
    @(19.12245521472685, 9.164981928519584, 5.960965703723132, 9.925660138972253, 43.2476179279168, 17.42058733001563, 26.11172261451397, 10.265873657791, 18.18778313370869),
    @(20.40225678621007, 10.48361110621365, 6.250484688405425, 9.872623844106432, 44.31443396404969, 17.43265455583065, 26.06891077066033, 10.33503125275632, 18.01299265805966),
    @(21.31398396833786, 11.35015019409668, 6.463522719436524, 9.838259996086096, 45.16115397479683, 17.46311476956634, 26.07222607222093, 10.39574410846939, 17.89380549987717),
    @(21.72049507559245, 11.72222360037, 6.56002358019707, 9.823620891542054, 45.55848204871284, 17.48166850972155, 26.08134777770834, 10.42545594917729, 17.84156224214298),
    @(21.87309270712595, 11.85994500524912, 6.596468442068196, 9.818219864357502, 45.71055211143803, 17.48936915793805, 26.08590037100979, 10.43700239769796, 17.82206118231936),
    @(21.84028995942536, 11.83042531730981, 6.588624340888733, 9.819376740445037, 45.67773168131406, 17.48768069472813, 26.08487099016437, 10.43450263008067, 17.82624855576657),
    @(21.73307705328237, 11.73361760815911, 6.563024098641863, 9.823173692086055, 45.5709614603149, 17.48228853029382, 26.08170029164276, 10.42640000045778, 17.83995222999031),
    @(21.66722728021995, 11.67390686003362, 6.54732934089342, 9.825517980407307, 45.50576724790772, 17.47907351299626, 26.0799012812494, 10.42147517315258, 17.84838284151089),
    @(21.28723820286, 11.32538969324265, 6.457204527263194, 9.839236652633538, 45.13542054131735, 17.46199671523549, 26.07178347975282, 10.39384400633955, 17.8972593229267),
    @(21.05190279545765, 11.10592403334839, 6.401782354320451, 9.847906751558035, 44.91124072836731, 17.45272345086997, 26.06875634028938, 10.37742533060575, 17.92774820152562),
    @(20.91577704323521, 10.9776106929566, 6.369868486670106, 9.852987065780773, 44.78345055662752, 17.44783192834069, 26.06773168067534, 10.36817914397952, 17.94547065363662),
    @(20.86956039945792, 10.93380817703729, 6.359057943651409, 9.854723242339498, 44.74038499228048, 17.44625170331159, 26.06750767807802, 10.36508261734647, 17.95150318024247),
    @(21.07703520227395, 11.12950201995512, 6.407686203853086, 9.84697413000308, 44.93498679230201, 17.4536648412877, 26.06900440369143, 10.37915273698385, 17.9244833651737),
    @(21.76460559549047, 11.76213848263619, 6.570546467357097, 9.822054571009277, 45.60227985703537, 17.4838540363589, 26.08260176862296, 10.42877197523207, 17.83591948036228),
    @(22.2061039904867, 12.15710231834522, 6.676397723567479, 9.806598553731597, 46.04771667493495, 17.5075165451304, 26.09789172955278, 10.4629188036492, 17.7796826744127),
    @(21.97123553009389, 11.94799304119555, 6.619968963214347, 9.814771852328841, 45.80917155876029, 17.49452804717174, 26.08914435578737, 10.44453884048837, 17.80954739484406),
    @(21.06567540917471, 11.11884907717264, 6.405017229150866, 9.84739546991521, 44.92424777938648, 17.45323786865509, 26.06889002517939, 10.37837117456365, 17.92595879326278),
    @(20.06030425080431, 10.14485763294566, 6.171934107803475, 9.88616159039433, 44.01427772156998, 17.42560241615288, 26.07441085495081, 10.31456481404629, 18.05864812069107)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $colLetters.Count; $j++) {
        $ws.Range($colLetters[$j] + $row).Value = $rowValues[$j]
    }
}
